$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")
$vals = @(0.9999975483099176, 0.9990478738077048, 0.9999954602780713, 0.9999872206564023, 0.9999927072274606, 0.000002288545075197367, 0.0008887680069980539, 0.000002735029408072741, 0.00000518055027330658, 0.00000395778984068966, 0.00009506688379990139, 0.001512793797976898, 1.000002028984896, 0.001577196561871648, 131.9751885613056, 196.5756072793203)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
